$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("TOLENTINO VASQUEZ DIANA KATHERYN", 15),
    @("JULCA VALENZUELA CINTIA KARYN", 15),
    @("LEON VERA MELISSA FIORELLA", 13),
    @("CARRILLO MARTÍNEZ HEIDY NAYELI", 12),
    @("CYNTHIA RODRIGUEZ LECCA", 12),
    @("RUBIO MARIÑOS GISELA JUDITH", 12),
    @("RUTH MELISSA RAMIREZ VELEZMORO", 11),
    @("DE LA CRUZ BENITES RICHARD ALEXANDER", 11),
    @("VALER VEGA PATRICIA GERALDINE", 11),
    @("SANCHEZ CORTEZ LEYLA DIANA", 11),
    @("ARENAS ZAVALA ANDYELA PATRICIA ISIDORA", 10),
    @("GASLAC GUTIERREZ FRANK JHORDY", 10),
    @("REYES RODRIGUEZ JEISSON STEVEN", 10),
    @("YZQUIERDO CARHUATANTA LEYDY YANELA", 10),
    @("PIERINA NAGIELLY SANDOVAL CONTRERAS", 9),
    @("PONCE VILLANUEVA CARMEN ISABEL", 9),
    @("SEGURA ASTO YAMILET ANTONELA", 9),
    @("GUZMAN ZAVALETA CECILIA MARISOL", 8),
    @("RODRIGUEZ RUBIO SANDRA MABEL", 8),
    @("RODRIGUEZ VASQUEZ WALTER", 2)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
